# Project "Sample Project" / Main.xlsx, sheet "Rules":
# Cell B11 currently holds the text "R40" (row for the 4th rule, "Good
# Night"). The commit changes that cell's content to the text "1" while
# keeping the cell's existing style (s="23") untouched.
#
# A plain  Range.Value = "1"  (or Value2/Formula) would be auto-typed by
# Excel as a *number*, which also silently rewrites the cell's style to a
# general/number format. To force the literal to remain a *text* value
# (so it round-trips as a shared string, like the original "R40"), we
# write it as a text formula and then collapse the formula down to a
# static value via copy / paste-special-values. That keeps the cell's
# type as text without touching its number format or style index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = $ws.Range("B11")
$target.Formula = '="1"'
$target.Copy($target) | Out-Null
$target.PasteSpecial(-4163) | Out-Null  # xlPasteValues
